{"js": "// Replace the date line and each \"AxB=\" multiplication prompt in the table\n// with its new value, preserving run formatting by using Range.insertText\n// with the \"Replace\" option (keeps the existing run's rPr).\nconst replacements = [\n  [\"2025-05-11 Sunday\", \"2025-05-12 Monday\"],\n  [\"70\u00d716=\", \"26\u00d733=\"],\n  [\"62\u00d744=\", \"11\u00d777=\"],\n  [\"60\u00d744=\", \"26\u00d725=\"],\n  [\"32\u00d791=\", \"76\u00d736=\"],\n  [\"92\u00d764=\", \"55\u00d768=\"],\n  [\"55\u00d735=\", \"65\u00d732=\"],\n  [\"53\u00d788=\", \"11\u00d721=\"],\n  [\"66\u00d764=\", \"55\u00d730=\"],\n  [\"54\u00d760=\", \"23\u00d773=\"],\n  [\"38\u00d772=\", \"11\u00d749=\"],\n  [\"42\u00d779=\", \"77\u00d775=\"],\n  [\"24\u00d727=\", \"71\u00d741=\"],\n  [\"12\u00d750=\", \"25\u00d777=\"],\n  [\"95\u00d772=\", \"39\u00d727=\"],\n  [\"45\u00d775=\", \"32\u00d752=\"],\n  [\"54\u00d757=\", \"56\u00d792=\"],\n  [\"45\u00d791=\", \"13\u00d732=\"],\n  [\"28\u00d778=\", \"32\u00d780=\"],\n  [\"69\u00d732=\", \"31\u00d751=\"],\n  [\"36\u00d730=\", \"47\u00d769=\"],\n  [\"58\u00d738=\", \"53\u00d798=\"],\n  [\"90\u00d755=\", \"26\u00d786=\"],\n  [\"57\u00d793=\", \"91\u00d786=\"],\n  [\"12\u00d734=\", \"82\u00d773=\"],\n  [\"63\u00d786=\", \"82\u00d782=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  // Only the first match should exist since every source string is unique,\n  // but guard against accidental repeats by replacing all occurrences found.\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"AxB=\" multiplication prompt in the table\n# with its new value, using Find/Replace on Document.Content so the run's\n# existing character formatting (font, size) is preserved.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-11 Sunday\", \"2025-05-12 Monday\"),\n    @(\"70\u00d716=\", \"26\u00d733=\"),\n    @(\"62\u00d744=\", \"11\u00d777=\"),\n    @(\"60\u00d744=\", \"26\u00d725=\"),\n    @(\"32\u00d791=\", \"76\u00d736=\"),\n    @(\"92\u00d764=\", \"55\u00d768=\"),\n    @(\"55\u00d735=\", \"65\u00d732=\"),\n    @(\"53\u00d788=\", \"11\u00d721=\"),\n    @(\"66\u00d764=\", \"55\u00d730=\"),\n    @(\"54\u00d760=\", \"23\u00d773=\"),\n    @(\"38\u00d772=\", \"11\u00d749=\"),\n    @(\"42\u00d779=\", \"77\u00d775=\"),\n    @(\"24\u00d727=\", \"71\u00d741=\"),\n    @(\"12\u00d750=\", \"25\u00d777=\"),\n    @(\"95\u00d772=\", \"39\u00d727=\"),\n    @(\"45\u00d775=\", \"32\u00d752=\"),\n    @(\"54\u00d757=\", \"56\u00d792=\"),\n    @(\"45\u00d791=\", \"13\u00d732=\"),\n    @(\"28\u00d778=\", \"32\u00d780=\"),\n    @(\"69\u00d732=\", \"31\u00d751=\"),\n    @(\"36\u00d730=\", \"47\u00d769=\"),\n    @(\"58\u00d738=\", \"53\u00d798=\"),\n    @(\"90\u00d755=\", \"26\u00d786=\"),\n    @(\"57\u00d793=\", \"91\u00d786=\"),\n    @(\"12\u00d734=\", \"82\u00d773=\"),\n    @(\"63\u00d786=\", \"82\u00d782=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
